$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows
$ws.Range("B2").Value = 10
$ws.Range("B3").Value = 7.5

# Add two new rows of PV component data
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 5
$ws.Range("C4").Value = "kW_peak"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = "kW_peak"

[void]$ws.Range("C8").Select()
$excel.ActiveWindow.Zoom = 228
